$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 441, shifting existing rows (441..539) down to (442..540)
$ws.Rows(441).Insert()

# Populate the newly inserted row 441 with the new record.
# (Same market/category context as neighboring rows; differing fields per the
# commit: Fecha, Variedad, Volumen, Unidad de comercializacion, Origen,
# Precio $/Kg and Kg / unidad.)
$ws.Cells.Item(441, 1).Value = 5
$ws.Cells.Item(441, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(441, 3).Value = "Maule"
$ws.Cells.Item(441, 4).Value = 45173
$ws.Cells.Item(441, 5).Value = 7
$ws.Cells.Item(441, 6).Value = "Fruta"
$ws.Cells.Item(441, 7).Value = 100102
$ws.Cells.Item(441, 8).Value = "Cítricos"
$ws.Cells.Item(441, 9).Value = 100102004
$ws.Cells.Item(441, 10).Value = "Mandarina"
$ws.Cells.Item(441, 11).Value = "Murcott"
$ws.Cells.Item(441, 12).Value = "Primera"
$ws.Cells.Item(441, 13).Value = 300
$ws.Cells.Item(441, 14).Value = 8000
$ws.Cells.Item(441, 15).Value = 8000
$ws.Cells.Item(441, 16).Value = 8000
$ws.Cells.Item(441, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(441, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(441, 19).Value = 444
$ws.Cells.Item(441, 20).Value = 18

# Ensure the date cell keeps the same date number format used by the rest of
# column D (style index 2, format "YYYY-MM-DD HH:MM:SS").
$ws.Cells.Item(441, 4).NumberFormat = $ws.Cells.Item(442, 4).NumberFormat
